$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.652.83'
$ws.Range("E2").Value = '  +0.19%  '

$ws.Range("D3").Value = '1.949.74'
$ws.Range("E3").Value = '  +1.28%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9970'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.27%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.72'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.29%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9969'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.32%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4846'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +2.45%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2920'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.30%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06825'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.27%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '112.51'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +6.14%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.64'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +6.60%  '

$ws.Range("D12").Value = '1.934.11'
$ws.Range("E12").Value = '  +0.48%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.490'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +3.04%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.07584'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.80%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6850'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.92%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '299.43'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.05%  '

$ws.Range("D17").Value = '30.584.21'
$ws.Range("E17").Value = '  -0.11%  '

$ws.Range("B18").Value = 'Avalanche'
$ws.Range("C18").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.16'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.54%  '

$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007699'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.67%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.600'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.22%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9972'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.24%  '

$ws.Range("D22").Value = '2.177.95'
$ws.Range("E22").Value = '  +0.31%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9965'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.37%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.527'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.58%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.541'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.04%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.73'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.08%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.63'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.47%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.152'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.95%  '

$ws.Range("E29").Value = '  +0.20%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.449'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +2.99%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.179'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.41%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.108'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +1.02%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04999'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.72%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7451'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.34%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.157'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.82%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02045'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.95%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.709'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.81%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.706'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.72%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.046'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.12%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '110.42'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.07%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4490'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.73%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8720'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.20%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.888'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.19%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '70.01'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +2.93%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9995'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.05%  '

$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.318'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.05%  '

$ws.Range("B47").Value = 'BitcoinSV'
$ws.Range("C47").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '49.48'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.92%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.333'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.51%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1240'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.10%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.2546'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.65%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '35.16'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.16%  '
